$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 2021 data row (row 10) following the same layout as existing rows.
$ws.Range("A10").Value = "2021年"

$values = @(41709, 1245482, 62253, 26086, 400437, 6960, 9179, 136752, 122911, 309585, 5240, 13651, 5809, 16918, 27433, 49863, 10696)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2  # Column B is index 2
    $ws.Cells.Item(10, $col).Value = $values[$i]
}

# Copy the row-9 cell format (bold/border/alignment) onto the new A10 label
# cell so it matches the styling of the other year-label cells (A2:A9).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
